# Daily attendance processing - 2026-01-29 18:58:50
# Swap the order of the two names recorded in the "Recorded By" column (G)
# from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$dim = $ws.UsedRange
$lastRow = $dim.Row + $dim.Rows.Count - 1

$searchRange = $ws.Range("G1:G" + $lastRow)

$first = $searchRange.Find($oldValue)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $current = $first
    do {
        $current.Value2 = $newValue
        $current = $searchRange.FindNext($current)
    } while (($current -ne $null) -and ($current.Address() -ne $firstAddress))
}
